# Students.xlsx edit — "Social Capital 2 done"
#
# 1. Two students previously marked against the (now retired) "6.2 Anderson"
#    presentation slot are reassigned to "7.2 Anderson".
# 2. The sheet is re-sorted by the Presentation column (column C), ascending,
#    with blank presentations sorting to the bottom (Excel default for Sort).
# 3. A new student, "Zoe", is appended with a "5.1 Calarco" presentation.
# 4. A previously-unassigned student (Ellermann, Kristine Grosen) is assigned
#    to "7.2 Anderson" after the sort.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: reassign "6.2 Anderson" -> "7.2 Anderson" for the two students
#     currently on that row (before the sort moves things around). ---
$ws.Range("C8").Value = "7.2 Anderson"
$ws.Range("C22").Value = "7.2 Anderson"

# --- Step 2: sort the data range A2:C49 by column C ascending ---
$dataRange = $ws.Range("A2:C49")
$sortKey = $ws.Range("C2:C49")
$dataRange.Sort($sortKey, 1)

# --- Step 3: append the new student "Zoe" in row 50 ---
$ws.Range("A50").Value = "Zoe"
$ws.Range("C50").Value = "5.1 Calarco"

# --- Step 4: assign Ellermann, Kristine Grosen (now on row 46 after the
#     sort pushed the blank-presentation rows to the bottom) to "7.2 Anderson" ---
$ws.Range("C46").Value = "7.2 Anderson"

# --- Update the view: scroll so row 16 is at the top and select F49,
#     matching the author's on-screen state when they saved. ---
$ws.Application.ActiveWindow.ScrollRow = 16
$ws.Range("F49").Select()

Write-Output "done"
